$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> (DAMSLTag, DialogAct)
$updates = @{
    42  = @("aa", "Agree/Accept")
    59  = @("sd", "Statement-non-opinion")
    60  = @("sv", "Statement-opinion")
    81  = @("aa", "Agree/Accept")
    86  = @("%", "Uninterpretable")
    89  = @("sd", "Statement-non-opinion")
    90  = @("sd", "Statement-non-opinion")
    103 = @("%", "Uninterpretable")
    105 = @("sd", "Statement-non-opinion")
    108 = @("sd", "Statement-non-opinion")
    110 = @("sd", "Statement-non-opinion")
    118 = @("ba", "Appreciation")
    119 = @("sv", "Statement-opinion")
    152 = @("qy", "Yes-No-Question")
    158 = @("sv", "Statement-opinion")
    164 = @("sv", "Statement-opinion")
    166 = @("sd", "Statement-non-opinion")
    167 = @("aa", "Agree/Accept")
    170 = @("sv", "Statement-opinion")
    174 = @("sv", "Statement-opinion")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
